# Apply the pilotA.xlsx view-state edit:
#  - PilotAReferenceSheet (sheet 3) zooms out to 45% and loses the tab
#    selection / loses being "on top".
#  - ListGeneration (sheet 1) becomes the active/selected tab, scrolled so
#    column S is visible, with the whole S column selected.
#  - RAND()-driven helper columns (B) on ListGeneration and
#    PilotAListGeneration recalc naturally as part of the volatile-function
#    refresh that happens whenever the workbook is recalculated/saved.

$wb = $excel.ActiveWorkbook

$wsList        = $wb.Worksheets.Item("ListGeneration")
$wsPilotList   = $wb.Worksheets.Item("PilotAListGeneration")
$wsReference   = $wb.Worksheets.Item("PilotAReferenceSheet")
$wsData        = $wb.Worksheets.Item("PilotAData")

# Zoom the reference sheet out to 45% (it was the tabbed/active sheet
# before the edit) before moving the active selection away from it.
$wsReference.Activate()
$excel.ActiveWindow.Zoom = 45

# Make ListGeneration the active sheet/tab and select the S column
# (S1:S1048576), matching the new selection left behind in the sheet.
$wsList.Activate()
$wsList.Range("S1:S1048576").Select()

# Force a full recalculation so every RAND() cell (ListGeneration!B1:B171
# and PilotAListGeneration!B1:B21) gets a fresh volatile value, same as
# what happens whenever the workbook is opened/recalculated in Excel.
$excel.CalculateFull()
